# ГОСТЬ.xlsx update
# - Rename sheet "М-10-2, А-10-2" -> "ГОСТЬ"
# - Add a new shared string "DEMO" as the value of A1
# - Swap the box-border formatting between G2 and G13 (G2 loses its
#   border/becomes plain, G13 gains the border the other day-blocks use)
# - Move the active selection from G12 to G3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "ГОСТЬ"

# Put the new "DEMO" label in A1 (keeps A1's existing style s="1")
$ws.Range("A1").Value = "DEMO"

# G2 should pick up the plain style currently used by G1/G14 (no border)
$ws.Range("G1").Copy()
$ws.Range("G2").PasteSpecial(-4122)

# G13 should pick up the bordered "box" style used by the rest of the
# Monday block (e.g. G3..G12)
$ws.Range("G3").Copy()
$ws.Range("G13").PasteSpecial(-4122)

# Clear clipboard / marching ants
$excel.CutCopyMode = 0

# Update the saved cursor/selection to G3
$ws.Range("G3").Select()

Write-Host "edit.ps1 applied"
